$d = $word.ActiveDocument

$d.Content.Find.Execute('94×15=1410', $true, $false, $false, $false, $false, $true, 1, $false, '85×57=4845', 2) | Out-Null
$d.Content.Find.Execute('25×35=875', $true, $false, $false, $false, $false, $true, 1, $false, '92×86=7912', 2) | Out-Null
$d.Content.Find.Execute('83×36=2988', $true, $false, $false, $false, $false, $true, 1, $false, '36×98=3528', 2) | Out-Null
$d.Content.Find.Execute('89×36=3204', $true, $false, $false, $false, $false, $true, 1, $false, '83×95=7885', 2) | Out-Null
$d.Content.Find.Execute('78×67=5226', $true, $false, $false, $false, $false, $true, 1, $false, '49×73=3577', 2) | Out-Null
$d.Content.Find.Execute('23×25=575', $true, $false, $false, $false, $false, $true, 1, $false, '71×40=2840', 2) | Out-Null
$d.Content.Find.Execute('15×20=300', $true, $false, $false, $false, $false, $true, 1, $false, '50×85=4250', 2) | Out-Null
$d.Content.Find.Execute('52×37=1924', $true, $false, $false, $false, $false, $true, 1, $false, '17×48=816', 2) | Out-Null
$d.Content.Find.Execute('56×56=3136', $true, $false, $false, $false, $false, $true, 1, $false, '33×65=2145', 2) | Out-Null
$d.Content.Find.Execute('96×80=7680', $true, $false, $false, $false, $false, $true, 1, $false, '26×94=2444', 2) | Out-Null
$d.Content.Find.Execute('47×95=4465', $true, $false, $false, $false, $false, $true, 1, $false, '67×15=1005', 2) | Out-Null
$d.Content.Find.Execute('28×89=2492', $true, $false, $false, $false, $false, $true, 1, $false, '34×13=442', 2) | Out-Null
$d.Content.Find.Execute('73×33=2409', $true, $false, $false, $false, $false, $true, 1, $false, '60×13=780', 2) | Out-Null
$d.Content.Find.Execute('25×52=1300', $true, $false, $false, $false, $false, $true, 1, $false, '68×89=6052', 2) | Out-Null
$d.Content.Find.Execute('27×59=1593', $true, $false, $false, $false, $false, $true, 1, $false, '51×55=2805', 2) | Out-Null
$d.Content.Find.Execute('83×81=6723', $true, $false, $false, $false, $false, $true, 1, $false, '69×69=4761', 2) | Out-Null
$d.Content.Find.Execute('59×65=3835', $true, $false, $false, $false, $false, $true, 1, $false, '54×18=972', 2) | Out-Null
$d.Content.Find.Execute('15×71=1065', $true, $false, $false, $false, $false, $true, 1, $false, '91×79=7189', 2) | Out-Null
$d.Content.Find.Execute('75×24=1800', $true, $false, $false, $false, $false, $true, 1, $false, '67×89=5963', 2) | Out-Null
$d.Content.Find.Execute('70×22=1540', $true, $false, $false, $false, $false, $true, 1, $false, '36×46=1656', 2) | Out-Null
$d.Content.Find.Execute('19×52=988', $true, $false, $false, $false, $false, $true, 1, $false, '34×81=2754', 2) | Out-Null
$d.Content.Find.Execute('67×22=1474', $true, $false, $false, $false, $false, $true, 1, $false, '27×37=999', 2) | Out-Null
$d.Content.Find.Execute('51×24=1224', $true, $false, $false, $false, $false, $true, 1, $false, '59×56=3304', 2) | Out-Null
$d.Content.Find.Execute('29×81=2349', $true, $false, $false, $false, $false, $true, 1, $false, '93×99=9207', 2) | Out-Null
$d.Content.Find.Execute('80×60=4800', $true, $false, $false, $false, $false, $true, 1, $false, '94×61=5734', 2) | Out-Null
